$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking strings
# (e.g. '1.002', '216.96') are not auto-converted to numbers/dates.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '26.152.93'
$ws.Range("E2").Value = '  -0.44%  '

# Row 3
$ws.Range("D3").Value = '1.671.93'
$ws.Range("E3").Value = '  -0.39%  '

# Row 4
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.18%  '

# Row 5
$ws.Range("D5").Value = '216.96'
$ws.Range("E5").Value = '  -0.87%  '

# Row 6
$ws.Range("E6").Value = '  +1.56%  '

# Row 7
$ws.Range("E7").Value = '  -0.16%  '

# Row 8
$ws.Range("D8").Value = '0.2704'
$ws.Range("E8").Value = '  +1.39%  '

# Row 9
$ws.Range("E9").Value = '  +0.40%  '

# Row 10
$ws.Range("E10").Value = '  -0.60%  '

# Row 11
$ws.Range("D11").Value = '0.07439'
$ws.Range("E11").Value = '  +0.78%  '

# Row 12
$ws.Range("D12").Value = '1.692.83'
$ws.Range("E12").Value = '  +0.83%  '

# Row 13
$ws.Range("D13").Value = '4.524'

# Row 14
$ws.Range("D14").Value = '0.5836'
$ws.Range("E14").Value = '  +0.98%  '

# Row 15
$ws.Range("D15").Value = '0.000008518'
$ws.Range("E15").Value = '  -1.20%  '

# Row 16
$ws.Range("D16").Value = '64.25'
$ws.Range("E16").Value = '  -1.31%  '

# Row 17
$ws.Range("D17").Value = '25.932.63'
$ws.Range("E17").Value = '  -1.65%  '

# Row 18
$ws.Range("D18").Value = '4.935'
$ws.Range("E18").Value = '  -1.76%  '

# Row 19
$ws.Range("E19").Value = '  -0.21%  '

# Row 20
$ws.Range("D20").Value = '10.81'
$ws.Range("E20").Value = '  -0.82%  '

# Row 21
$ws.Range("D21").Value = '190.02'
$ws.Range("E21").Value = '  +1.44%  '

# Row 22
$ws.Range("D22").Value = '6.197'
$ws.Range("E22").Value = '  -0.74%  '

# Row 23
$ws.Range("D23").Value = '1.003'

# Row 24
$ws.Range("D24").Value = '144.78'
$ws.Range("E24").Value = '  +0.76%  '

# Row 25
$ws.Range("D25").Value = '0.1242'
$ws.Range("E25").Value = '  +5.02%  '

# Row 26
$ws.Range("D26").Value = '7.619'
$ws.Range("E26").Value = '  +0.62%  '

# Row 27
$ws.Range("D27").Value = '15.73'
$ws.Range("E27").Value = '  -0.12%  '

# Row 28
$ws.Range("D28").Value = '0.06561'
$ws.Range("E28").Value = '  +12.52%  '

# Row 29
$ws.Range("E29").Value = '  +0.72%  '

# Row 30
$ws.Range("D30").Value = '1.317'
$ws.Range("E30").Value = '  -0.98%  '

# Row 31
$ws.Range("D31").Value = '3.590'
$ws.Range("E31").Value = '  +1.99%  '

# Row 32
$ws.Range("D32").Value = '3.530'
$ws.Range("E32").Value = '  +0.54%  '

# Row 33
$ws.Range("E33").Value = '  +0.05%  '

# Row 34
$ws.Range("D34").Value = '1.020'
$ws.Range("E34").Value = '  +1.33%  '

# Row 35
$ws.Range("D35").Value = '0.6176'
$ws.Range("E35").Value = '  +2.94%  '

# Row 36
$ws.Range("D36").Value = '2.366'
$ws.Range("E36").Value = '  +0.29%  '

# Row 37
$ws.Range("E37").Value = '  +1.21%  '

# Row 38
$ws.Range("D38").Value = '6.261'
$ws.Range("E38").Value = '  +6.20%  '

# Row 39
$ws.Range("D39").Value = '1.094.51'
$ws.Range("E39").Value = '  -0.60%  '

# Row 40
$ws.Range("D40").Value = '0.01599'
$ws.Range("E40").Value = '  -0.92%  '

# Row 41
$ws.Range("D41").Value = '0.8710'
$ws.Range("E41").Value = '  +0.97%  '

# Row 42
$ws.Range("D42").Value = '1.011'
$ws.Range("E42").Value = '  +0.65%  '

# Row 43
$ws.Range("D43").Value = '100.90'
$ws.Range("E43").Value = '  +1.33%  '

# Row 44
$ws.Range("D44").Value = '1.817.70'
$ws.Range("E44").Value = '  -0.60%  '

# Row 45
$ws.Range("E45").Value = '  -1.30%  '

# Row 46
$ws.Range("D46").Value = '56.53'
$ws.Range("E46").Value = '  +0.01%  '

# Row 47
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '8.155'
$ws.Range("E47").Value = '  +0.96%  '

# Row 48
$ws.Range("B48").Value = 'Frax'
$ws.Range("C48").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D48").Value = '1.005'
$ws.Range("E48").Value = '  +0.05%  '

# Row 49
$ws.Range("D49").Value = '0.05236'
$ws.Range("E49").Value = '  +0.26%  '

# Row 50
$ws.Range("D50").Value = '0.4278'
$ws.Range("E50").Value = '  -0.83%  '

# Row 51
$ws.Range("D51").Value = '6.000'
$ws.Range("E51").Value = '  +1.81%  '

# Restore default (general) formatting on column D so the saved
# style matches the original (no explicit style index).
$dRange.ClearFormats()
